# Updates cryptos list prices/volumes (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.871.53"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.31"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.57"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.46"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.641.34"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.61"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.856.30"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.76"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.96"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.24"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.78"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.42"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0496"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.130.22"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.46"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.41"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.798"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.763.29"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.16"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("E48").Value = "  +4.18%  "
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.57"
